$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = ""

$ws.Range("H18").Value = 1350
$ws.Range("I18").Value = 1350
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1350
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1066

$ws.Range("H33").Value = 820.75
$ws.Range("I33").Value = 94.5
$ws.Range("J33").Value = 2999.5
$ws.Range("K33").Value = 94.5
$ws.Range("L33").Value = 2999.5
$ws.Range("M33").Value = 134.5
$ws.Range("N33").Value = -3457.5

$ws.Range("H41").Value = 295.2857
$ws.Range("I41").Value = 233.6
$ws.Range("J41").Value = 449.5
$ws.Range("K41").Value = 233.6
$ws.Range("L41").Value = 449.5
$ws.Range("M41").Value = 206.4
$ws.Range("N41").Value = -1329.5

$ws.Range("H132").Value = 2530.5293
$ws.Range("I132").Value = 2530.5293
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7591.5879
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5061.5879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 778.2
$ws.Range("I2").Value = 597.75
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 597.75
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -484.75
$ws.Range("N2").Value = -1726

$ws.Range("H61").Value = 3589.3333
$ws.Range("I61").Value = 3926.4
$ws.Range("J61").Value = 1904
$ws.Range("K61").Value = 3926.4
$ws.Range("L61").Value = 1904
$ws.Range("M61").Value = -3714.4
$ws.Range("N61").Value = -2328

$ws.Range("H110").Value = 511.57144
$ws.Range("I110").Value = 431
$ws.Range("J110").Value = 995
$ws.Range("K110").Value = 431
$ws.Range("L110").Value = 995
$ws.Range("M110").Value = 1614
$ws.Range("N110").Value = -5085

$ws.Range("H116").Value = 778.2
$ws.Range("I116").Value = 597.75
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 597.75
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 1696.25
$ws.Range("N116").Value = -6088

$ws.Range("H136").Value = 3589.3333
$ws.Range("I136").Value = 3926.4
$ws.Range("J136").Value = 1904
$ws.Range("K136").Value = 11779.2
$ws.Range("L136").Value = 5712
$ws.Range("M136").Value = -9229.200000000001
$ws.Range("N136").Value = -10812

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 778.2
$ws.Range("I3").Value = 597.75
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 597.75
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -483.75
$ws.Range("N3").Value = -1728

$ws.Range("H94").Value = 1232.5555
$ws.Range("I94").Value = 1385.375
$ws.Range("J94").Value = 10
$ws.Range("K94").Value = 1385.375
$ws.Range("L94").Value = 10
$ws.Range("M94").Value = -934.375
$ws.Range("N94").Value = -912

$ws.Range("H134").Value = 6801.3184
$ws.Range("I134").Value = 6721.5
$ws.Range("J134").Value = 7599.5
$ws.Range("K134").Value = 20164.5
$ws.Range("L134").Value = 22798.5
$ws.Range("M134").Value = -17629.5
$ws.Range("N134").Value = -27868.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1716.9524
$ws.Range("I31").Value = 1540.5555
$ws.Range("J31").Value = 2775.3333
$ws.Range("K31").Value = 1540.5555
$ws.Range("L31").Value = 2775.3333
$ws.Range("M31").Value = -1245.5555
$ws.Range("N31").Value = -3365.3333

$ws.Range("H34").Value = 1716.9524
$ws.Range("I34").Value = 1540.5555
$ws.Range("J34").Value = 2775.3333
$ws.Range("K34").Value = 1540.5555
$ws.Range("L34").Value = 2775.3333
$ws.Range("M34").Value = -1338.5555
$ws.Range("N34").Value = -3179.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 469.23077
$ws.Range("I9").Value = 442.85715
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 1328.57145
$ws.Range("L9").Value = 1500
$ws.Range("M9").Value = -1104.57145
$ws.Range("N9").Value = -1948

$ws.Range("H10").Value = 324
$ws.Range("I10").Value = 324
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 972
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -833

$ws.Range("H12").Value = 107.22222
$ws.Range("I12").Value = 30
$ws.Range("J12").Value = 129.28572
$ws.Range("K12").Value = 90
$ws.Range("L12").Value = 387.85716
$ws.Range("M12").Value = 83
$ws.Range("N12").Value = -733.85716

$ws.Range("H68").Value = 1799.6666
$ws.Range("I68").Value = 1799.6666
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5398.9998
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4587.9998
$ws.Range("N68").Value = ""

$ws.Range("H71").Value = 1799.6666
$ws.Range("I71").Value = 1799.6666
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 16196.9994
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -12140.9994
$ws.Range("N71").Value = ""

$ws.Range("H103").Value = 402.2
$ws.Range("I103").Value = 402.2
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1206.6
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -327.5999999999999
$ws.Range("N103").Value = ""

$ws.Range("H129").Value = 1256.25

$ws.Range("H137").Value = 1691.5
$ws.Range("I137").Value = 1350
$ws.Range("J137").Value = 2033
$ws.Range("K137").Value = 4050
$ws.Range("L137").Value = 6099
$ws.Range("M137").Value = 1050
$ws.Range("N137").Value = -16299

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 287
$ws.Range("I107").Value = 287
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 287
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1633

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 39813.145
$ws.Range("I7").Value = 39813.145
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 39813.145
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -39701.145

$ws.Range("H16").Value = 413.57144
$ws.Range("I16").Value = 448.83334
$ws.Range("J16").Value = 202
$ws.Range("K16").Value = 448.83334
$ws.Range("L16").Value = 202
$ws.Range("M16").Value = -278.83334
$ws.Range("N16").Value = -542

$ws.Range("H47").Value = 14874.25
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 14874.25
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 14874.25
$ws.Range("N47").Value = -15854.25

$ws.Range("H52").Value = 14874.25
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 14874.25
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 14874.25
$ws.Range("N52").Value = -15340.25

$ws.Range("H126").Value = 39813.145
$ws.Range("I126").Value = 39813.145
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 119439.435
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -116969.435

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1383.2222
$ws.Range("I4").Value = 287
$ws.Range("J4").Value = 2260.2
$ws.Range("K4").Value = 287
$ws.Range("L4").Value = 2260.2
$ws.Range("M4").Value = -174
$ws.Range("N4").Value = -2486.2

$ws.Range("H92").Value = 33332
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 33332
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 33332
$ws.Range("N92").Value = -38324

$ws.Range("H100").Value = 945.2857
$ws.Range("I100").Value = 1011.6667
$ws.Range("J100").Value = 547
$ws.Range("K100").Value = 2023.3334
$ws.Range("L100").Value = 1094
$ws.Range("M100").Value = -1482.3334
$ws.Range("N100").Value = -2176
